$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Cong thuc san xuat" feature: add a "Loai san pham" column right before
# the existing "Danh muc" column (new column E), shifting old E:K to F:L.
$ws.Range("E1:E2").EntireColumn.Insert()

# Header text for the newly inserted column.
$ws.Range("E1").Value = "Loại sản phẩm (SP_NHA_CUNG_CAP, SP_SAN_XUAT, NGUYEN_LIEU)"

# Sample row: the new cell is left blank (it already inherited the centered
# style used by its D2/F2 neighbours from the column insert).

# The new column has no width yet - size it to roughly fit the long header.
$ws.Columns.Item(5).ColumnWidth = 59.8333333333333

# Restore the plain top-left view (no more horizontal scroll) and move the
# active selection to the newly-edited cell.
[void]$ws.Range("E9").Select()
